# Cronograma de Actividades - update schedule dates (documentation refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Diseño de la base de datos ---
# H3 (Comienzo) stays "25/09/17"
$ws.Range("J3").Value = "26/09/17"   # Fin: 27/09/17 -> 26/09/17

# --- Row 4: Diseño del diagrama de clases ---
$ws.Range("H4").Value = "27/09/17"      # Comienzo: 29/09/17 -> 27/09/17
$ws.Range("J4").Value = "28/09/2017"    # Fin: 10/2/2017 -> 28/09/2017 (text)

# --- Row 5: Implementar el diagrama de la base de datos ---
$ws.Range("H5").Value = "29/09/2017"    # Comienzo: 4/10/2017 -> 29/09/2017 (text)
$ws.Range("J5").Value = "29/09/2017"    # Fin: 4/10/2017 -> 29/09/2017 (text)

# --- Row 6: Implementacion de los diagramas de clases ---
$ws.Range("H6").Value = 42776
$ws.Range("H6").NumberFormat = "m/d/yy"
$ws.Range("J6").Value = 42776
$ws.Range("J6").NumberFormat = "m/d/yy"

# --- Row 7: Implementacion de los CRUDS necesarios del lado del servidor ---
$ws.Range("H7").Value = 42804
$ws.Range("H7").NumberFormat = "m/d/yy"
$ws.Range("J7").Value = 42988
$ws.Range("J7").NumberFormat = "m/d/yy"

# --- Row 8: Implementacion de los CRUDS necesarios del lado del Cliente ---
$ws.Range("H8").Value = 43018
$ws.Range("H8").NumberFormat = "m/d/yy"
$ws.Range("J8").Value = "16/10/17"      # Fin: 30/10/17 -> 16/10/17 (text)

# --- Row 9: Sincronizar el servidor con el cliente ---
$ws.Range("H9").Value = "17/10/2017"    # Comienzo -> text
$ws.Range("J9").Value = "18/10/2017"    # Fin -> text

# --- Row 10: Corregir Bugs si hubiera ---
$ws.Range("H10").Value = "19/10/2017"   # Comienzo -> text
$ws.Range("J10").Value = "20/10/2017"   # Fin -> text

# Selection marker moved from K18 to K15
$ws.Range("K15").Select()
